$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the series: insert a row at
# position 35 (pushing the existing rows 35-145 down to 36-146) and
# populate it with the new observation.
$ws.Rows("35:35").Insert()

$ws.Range("A35").Value = 5
$ws.Range("B35").Value = "Macroferia Regional de Talca"
$ws.Range("C35").Value = "Maule"
$ws.Range("D35").Value = 44811
$ws.Range("E35").Value = 7
$ws.Range("F35").Value = "Fruta"
$ws.Range("G35").Value = 100108
$ws.Range("H35").Value = "Tropicales y subtropicales"
$ws.Range("I35").Value = 100108002
$ws.Range("J35").Value = "Mango"
$ws.Range("K35").Value = "Sin especificar"
$ws.Range("L35").Value = "Primera"
$ws.Range("M35").Value = 228
$ws.Range("N35").Value = 9000
$ws.Range("O35").Value = 9000
$ws.Range("P35").Value = 9000
$ws.Range("Q35").Value = "$/bandeja 4 kilos"
$ws.Range("R35").Value = "Brasil"
$ws.Range("S35").Value = 2250
$ws.Range("T35").Value = 4
